$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.002.80"
$ws.Range("E2").Value = "  +1.39%  "
$ws.Range("D3").Value = "2.244.13"
$ws.Range("E3").Value = "  +0.37%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'317.84"
$ws.Range("E5").Value = "  +0.27%  "
$ws.Range("D6").Value = "'100.57"
$ws.Range("E6").Value = "  +1.11%  "
$ws.Range("E7").Value = "  -1.64%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  -3.59%  "
$ws.Range("D10").Value = "'36.76"
$ws.Range("E10").Value = "  -0.77%  "
$ws.Range("D11").Value = "'0.0826"
$ws.Range("E11").Value = "  -0.60%  "
$ws.Range("E12").Value = "  -2.47%  "
$ws.Range("E13").Value = "  -2.07%  "
$ws.Range("D14").Value = "2.584.78"
$ws.Range("E14").Value = "  +0.32%  "
$ws.Range("E15").Value = "  -2.21%  "
$ws.Range("D16").Value = "2.247.25"
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("D17").Value = "'14.12"
$ws.Range("E17").Value = "  -1.49%  "
$ws.Range("D18").Value = "43.905.84"
$ws.Range("E18").Value = "  +1.33%  "
$ws.Range("D19").Value = "'13.28"
$ws.Range("E19").Value = "  -6.32%  "
$ws.Range("D20").Value = "0.0₃0971"
$ws.Range("E20").Value = "  -0.20%  "
$ws.Range("E21").Value = "  -3.48%  "
$ws.Range("D22").Value = "'65.46"
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("D23").Value = "'3.08"
$ws.Range("E23").Value = "  -3.36%  "
$ws.Range("D24").Value = "'234.59"
$ws.Range("E24").Value = "  -0.68%  "
$ws.Range("D25").Value = "'2.05"
$ws.Range("E25").Value = "  -5.43%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").Value = "'10.44"
$ws.Range("E27").Value = "  +3.94%  "
$ws.Range("E28").Value = "  -0.54%  "
$ws.Range("D29").Value = "'37.44"
$ws.Range("E29").Value = "  +2.35%  "
$ws.Range("D30").Value = "'6.13"
$ws.Range("E30").Value = "  -4.01%  "
$ws.Range("D31").Value = "'159.55"
$ws.Range("E31").Value = "  +1.17%  "
$ws.Range("E32").Value = "  -1.26%  "
$ws.Range("D33").Value = "'0.0845"
$ws.Range("E33").Value = "  -2.99%  "
$ws.Range("E34").Value = "  -1.47%  "
$ws.Range("D35").Value = "'3.18"
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("E36").Value = "  +8.29%  "
$ws.Range("E37").Value = "  +2.67%  "
$ws.Range("E38").Value = "  -1.98%  "
$ws.Range("D39").Value = "'16.15"
$ws.Range("E39").Value = "  +12.13%  "
$ws.Range("D40").Value = "'3.65"
$ws.Range("E40").Value = "  -1.38%  "
$ws.Range("D41").Value = "'4.11"
$ws.Range("E41").Value = "  -6.07%  "
$ws.Range("D42").Value = "'0.0313"
$ws.Range("E42").Value = "  -2.47%  "
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("D44").Value = "1.743.90"
$ws.Range("E44").Value = "  -4.78%  "
$ws.Range("E45").Value = "  -2.85%  "
$ws.Range("D46").Value = "'81.62"
$ws.Range("E46").Value = "  -2.81%  "
$ws.Range("D47").Value = "'74.19"
$ws.Range("E47").Value = "  +0.47%  "
$ws.Range("D48").Value = "'5.14"
$ws.Range("E48").Value = "  -3.05%  "
$ws.Range("D49").Value = "'102.11"
$ws.Range("E49").Value = "  -1.04%  "
$ws.Range("E50").Value = "  +2.98%  "
$ws.Range("D51").Value = "'57.33"
$ws.Range("E51").Value = "  -1.46%  "
